$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2235.4443
$ws.Range("I4").Value = 2235.4443
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 2235.4443
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -2121.4443
$ws.Range("N4").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29590.355
$ws.Range("I32").Value = 4762.135
$ws.Range("K32").Value = 4762.135
$ws.Range("M32").Value = -4475.135
$ws.Range("H61").Value = 1469.6316
$ws.Range("I61").Value = 1288.7646
$ws.Range("J61").Value = 3007
$ws.Range("K61").Value = 1288.7646
$ws.Range("L61").Value = 3007
$ws.Range("M61").Value = -1076.7646
$ws.Range("N61").Value = -3431
$ws.Range("H106").Value = 43990
$ws.Range("J106").Value = 43990
$ws.Range("L106").Value = 43990
$ws.Range("N106").Value = -46514
$ws.Range("H132").Value = 1951.2413
$ws.Range("I132").Value = 1629.619
$ws.Range("J132").Value = 2795.5
$ws.Range("K132").Value = 4888.857
$ws.Range("L132").Value = 8386.5
$ws.Range("M132").Value = -2358.857
$ws.Range("N132").Value = -13446.5
$ws.Range("H136").Value = 1469.6316
$ws.Range("I136").Value = 1288.7646
$ws.Range("J136").Value = 3007
$ws.Range("K136").Value = 3866.2938
$ws.Range("L136").Value = 9021
$ws.Range("M136").Value = -1316.2938
$ws.Range("N136").Value = -14121

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 54334.39
$ws.Range("I86").Value = 72930.06
$ws.Range("J86").Value = 1646.6666
$ws.Range("K86").Value = 72930.06
$ws.Range("L86").Value = 1646.6666
$ws.Range("M86").Value = -71807.06
$ws.Range("N86").Value = -3892.6666
$ws.Range("H89").Value = 54334.39
$ws.Range("I89").Value = 72930.06
$ws.Range("J89").Value = 1646.6666
$ws.Range("K89").Value = 364650.3
$ws.Range("L89").Value = 8233.333000000001
$ws.Range("M89").Value = -359034.3
$ws.Range("N89").Value = -19465.333

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13012.954
$ws.Range("I58").Value = 1303.7368
$ws.Range("J58").Value = 87171.336
$ws.Range("K58").Value = 1303.7368
$ws.Range("L58").Value = 87171.336
$ws.Range("M58").Value = -1100.7368
$ws.Range("N58").Value = -87577.336
$ws.Range("H59").Value = 21886.924
$ws.Range("J59").Value = 21886.924
$ws.Range("L59").Value = 21886.924
$ws.Range("N59").Value = -24176.924
$ws.Range("H122").Value = 4331.2383
$ws.Range("I122").Value = 3395.5625
$ws.Range("K122").Value = 10186.6875
$ws.Range("M122").Value = -7736.6875
$ws.Range("H132").Value = 48390740
$ws.Range("I132").Value = 47623050
$ws.Range("J132").Value = 50002892
$ws.Range("K132").Value = 142869150
$ws.Range("L132").Value = 150008676
$ws.Range("M132").Value = -142866620
$ws.Range("N132").Value = -150013736
$ws.Range("H134").Value = 2396.2
$ws.Range("I134").Value = 2495.25
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 7485.75
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = -4950.75
$ws.Range("N134").Value = -11070
$ws.Range("H136").Value = 13012.954
$ws.Range("I136").Value = 1303.7368
$ws.Range("J136").Value = 87171.336
$ws.Range("K136").Value = 3911.2104
$ws.Range("L136").Value = 261514.008
$ws.Range("M136").Value = -1361.2104
$ws.Range("N136").Value = -266614.008

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 50.5625
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 50.5625
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 151.6875
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -497.6875
$ws.Range("H122").Value = 6888.5
$ws.Range("J122").Value = 15199.571
$ws.Range("L122").Value = 136796.139
$ws.Range("N122").Value = -141696.139
$ws.Range("H131").Value = 857.2
$ws.Range("I131").Value = 538
$ws.Range("J131").Value = 874
$ws.Range("K131").Value = 1614
$ws.Range("L131").Value = 2622
$ws.Range("M131").Value = 3426
$ws.Range("N131").Value = -12702

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 111126056
$ws.Range("I97").Value = 111126056
$ws.Range("K97").Value = 111126056
$ws.Range("M97").Value = -111125560
$ws.Range("H102").Value = 2103.9722
$ws.Range("I102").Value = 2062.8965
$ws.Range("J102").Value = 2274.1428
$ws.Range("K102").Value = 2062.8965
$ws.Range("L102").Value = 2274.1428
$ws.Range("M102").Value = -440.8964999999998
$ws.Range("N102").Value = -5518.1428
$ws.Range("H113").Value = 2074.0715
$ws.Range("I113").Value = 2870.2
$ws.Range("K113").Value = 2870.2
$ws.Range("M113").Value = -700.1999999999998
$ws.Range("H126").Value = 2128.6667
$ws.Range("I126").Value = 1971.7142
$ws.Range("J126").Value = 2193.2942
$ws.Range("K126").Value = 5915.142599999999
$ws.Range("L126").Value = 6579.882599999999
$ws.Range("M126").Value = -3445.142599999999
$ws.Range("N126").Value = -11519.8826

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2110.889
$ws.Range("I7").Value = 1649.9445
$ws.Range("J7").Value = 3032.7778
$ws.Range("K7").Value = 1649.9445
$ws.Range("L7").Value = 3032.7778
$ws.Range("M7").Value = -1537.9445
$ws.Range("N7").Value = -3256.7778
$ws.Range("H46").Value = 20000
$ws.Range("J46").Value = 20000
$ws.Range("L46").Value = 20000
$ws.Range("N46").Value = -20376
$ws.Range("H122").Value = 2188.5557
$ws.Range("I122").Value = 2114.2856
$ws.Range("J122").Value = 2448.5
$ws.Range("K122").Value = 6342.8568
$ws.Range("L122").Value = 7345.5
$ws.Range("M122").Value = -3892.8568
$ws.Range("N122").Value = -12245.5
$ws.Range("H126").Value = 2110.889
$ws.Range("I126").Value = 1649.9445
$ws.Range("J126").Value = 3032.7778
$ws.Range("K126").Value = 4949.833500000001
$ws.Range("L126").Value = 9098.3334
$ws.Range("M126").Value = -2479.833500000001
$ws.Range("N126").Value = -14038.3334

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 498.85715
$ws.Range("I113").Value = 384.35715
$ws.Range("K113").Value = 1153.07145
$ws.Range("M113").Value = 1016.92855
$ws.Range("H122").Value = 3875.625
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5001
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 15003
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -19903
$ws.Range("H132").Value = 2808.8965
$ws.Range("I132").Value = 2633.0435
$ws.Range("K132").Value = 7899.130500000001
$ws.Range("M132").Value = -5369.130500000001
$ws.Range("H136").Value = 1179.5
$ws.Range("I136").Value = 490.06897
$ws.Range("J136").Value = 2231.7896
$ws.Range("K136").Value = 1470.20691
$ws.Range("L136").Value = 6695.3688
$ws.Range("M136").Value = 1079.79309
$ws.Range("N136").Value = -11795.3688

